$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 82
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H9").Value = 568777.5600000001
$ws.Range("I9").Value = 1136619.9
$ws.Range("J9").Value = 935.1818
$ws.Range("K9").Value = 1136619.9
$ws.Range("L9").Value = 935.1818
$ws.Range("M9").Value = -1136450.9
$ws.Range("N9").Value = -1273.1818
$ws.Range("H18").Value = 4163.6665
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5568
$ws.Range("H28").Value = 3429.85
$ws.Range("I28").Value = 203.76923
$ws.Range("J28").Value = 9421.143
$ws.Range("K28").Value = 203.76923
$ws.Range("L28").Value = 9421.143
$ws.Range("M28").Value = 281.23077
$ws.Range("N28").Value = -10391.143
$ws.Range("H40").Value = 100002350
$ws.Range("J40").Value = 166668800
$ws.Range("L40").Value = 166668800
$ws.Range("N40").Value = -166669150
$ws.Range("H55").Value = 451.96155
$ws.Range("I55").Value = 143.79167
$ws.Range("K55").Value = 143.79167
$ws.Range("M55").Value = 70.20832999999999
$ws.Range("H112").Value = 3690.5833
$ws.Range("J112").Value = 4466.5557
$ws.Range("L112").Value = 13399.6671
$ws.Range("N112").Value = -15615.6671
$ws.Range("H129").Value = 3126.25
$ws.Range("I129").Value = 1072.625
$ws.Range("J129").Value = 4495.3335
$ws.Range("K129").Value = 3217.875
$ws.Range("L129").Value = 13486.0005
$ws.Range("M129").Value = 1782.125
$ws.Range("N129").Value = -23486.0005
$ws.Range("H137").Value = 31255456
$ws.Range("I137").Value = 83337140
$ws.Range("J137").Value = 6448.4
$ws.Range("K137").Value = 250011420
$ws.Range("L137").Value = 19345.2
$ws.Range("M137").Value = -250008870
$ws.Range("N137").Value = -24445.2
$ws.Range("H138").Value = 5903.881
$ws.Range("I138").Value = 3777.25
$ws.Range("J138").Value = 8739.388999999999
$ws.Range("K138").Value = 11331.75
$ws.Range("L138").Value = 26218.167
$ws.Range("M138").Value = -6191.75
$ws.Range("N138").Value = -36498.167

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1020
$ws.Range("I5").Value = 1020
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1020
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -908
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 5140.7676
$ws.Range("I32").Value = 4851.1025
$ws.Range("K32").Value = 4851.1025
$ws.Range("M32").Value = -4564.1025
$ws.Range("H45").Value = 3488.7036
$ws.Range("I45").Value = 1928.5333
$ws.Range("K45").Value = 1928.5333
$ws.Range("M45").Value = -1551.5333
$ws.Range("H74").Value = 2159.4827
$ws.Range("I74").Value = 1880.5
$ws.Range("K74").Value = 1880.5
$ws.Range("M74").Value = -1006.5
$ws.Range("H77").Value = 2159.4827
$ws.Range("I77").Value = 1880.5
$ws.Range("K77").Value = 9402.5
$ws.Range("M77").Value = -5034.5
$ws.Range("H132").Value = 2706714.2
$ws.Range("I132").Value = 4122.8887
$ws.Range("K132").Value = 12368.6661
$ws.Range("M132").Value = -9838.666100000002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1020
$ws.Range("I4").Value = 1020
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1020
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -905
$ws.Range("N4").ClearContents()
$ws.Range("H105").Value = 539639.8
$ws.Range("I105").Value = 992043.75
$ws.Range("J105").Value = 4980.5454
$ws.Range("K105").Value = 992043.75
$ws.Range("L105").Value = 4980.5454
$ws.Range("M105").Value = -990296.75
$ws.Range("N105").Value = -8474.545399999999
$ws.Range("H134").Value = 3573620.5
$ws.Range("I134").Value = 1731.4286
$ws.Range("K134").Value = 5194.2858
$ws.Range("M134").Value = -2659.2858

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23275958
$ws.Range("I31").Value = 50040212
$ws.Range("J31").Value = 2696.5217
$ws.Range("K31").Value = 50040212
$ws.Range("L31").Value = 2696.5217
$ws.Range("M31").Value = -50039917
$ws.Range("N31").Value = -3286.5217
$ws.Range("H33").Value = 6593.3335
$ws.Range("I33").Value = 1312
$ws.Range("J33").Value = 33000
$ws.Range("K33").Value = 1312
$ws.Range("L33").Value = 33000
$ws.Range("M33").Value = -933
$ws.Range("N33").Value = -33758
$ws.Range("H34").Value = 23275958
$ws.Range("I34").Value = 50040212
$ws.Range("J34").Value = 2696.5217
$ws.Range("K34").Value = 50040212
$ws.Range("L34").Value = 2696.5217
$ws.Range("M34").Value = -50040010
$ws.Range("N34").Value = -3100.5217
$ws.Range("H62").Value = 3884.3333
$ws.Range("I62").Value = 3560.1538
$ws.Range("J62").Value = 5991.5
$ws.Range("K62").Value = 3560.1538
$ws.Range("L62").Value = 5991.5
$ws.Range("M62").Value = -2936.1538
$ws.Range("N62").Value = -7239.5
$ws.Range("H65").Value = 3884.3333
$ws.Range("I65").Value = 3560.1538
$ws.Range("J65").Value = 5991.5
$ws.Range("K65").Value = 17800.769
$ws.Range("L65").Value = 29957.5
$ws.Range("M65").Value = -14680.769
$ws.Range("N65").Value = -36197.5
$ws.Range("H70").Value = 93332
$ws.Range("J70").Value = 93332
$ws.Range("L70").Value = 93332
$ws.Range("N70").Value = -93962
$ws.Range("H73").Value = 93332
$ws.Range("J73").Value = 93332
$ws.Range("L73").Value = 93332
$ws.Range("N73").Value = -95516
$ws.Range("H99").Value = 17073.867
$ws.Range("I99").Value = 5307.8335
$ws.Range("K99").Value = 5307.8335
$ws.Range("M99").Value = -3809.8335
$ws.Range("H126").Value = 17073.867
$ws.Range("I126").Value = 5307.8335
$ws.Range("K126").Value = 15923.5005
$ws.Range("M126").Value = -13453.5005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 7301.625
$ws.Range("J29").Value = 11666.375
$ws.Range("L29").Value = 34999.125
$ws.Range("N29").Value = -35553.125
$ws.Range("H56").Value = 15302.54
$ws.Range("I56").Value = 15302.54
$ws.Range("K56").Value = 15302.54
$ws.Range("M56").Value = -14772.54
$ws.Range("H80").Value = 41667936
$ws.Range("J80").Value = 1691
$ws.Range("L80").Value = 5073
$ws.Range("N80").Value = -6945
$ws.Range("H83").Value = 41667936
$ws.Range("J83").Value = 1691
$ws.Range("L83").Value = 15219
$ws.Range("N83").Value = -24579
$ws.Range("H104").Value = 19777.334
$ws.Range("I104").Value = 17999
$ws.Range("J104").Value = 20666.5
$ws.Range("K104").Value = 53997
$ws.Range("L104").Value = 61999.5
$ws.Range("M104").Value = -51376
$ws.Range("N104").Value = -67241.5
$ws.Range("H129").Value = 11908571
$ws.Range("I129").Value = 27781054
$ws.Range("J129").Value = 4209
$ws.Range("K129").Value = 83343162
$ws.Range("L129").Value = 12627
$ws.Range("M129").Value = -83338162
$ws.Range("N129").Value = -22627

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 969.4138
$ws.Range("I107").Value = 1004.75
$ws.Range("K107").Value = 1004.75
$ws.Range("M107").Value = 915.25
$ws.Range("H132").Value = 1820443.8
$ws.Range("I132").Value = 2260.8445
$ws.Range("J132").Value = 10002267
$ws.Range("K132").Value = 6782.5335
$ws.Range("L132").Value = 30006801
$ws.Range("M132").Value = -4252.5335
$ws.Range("N132").Value = -30011861

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2241.7742
$ws.Range("I40").Value = 2183.2
$ws.Range("K40").Value = 2183.2
$ws.Range("M40").Value = -2047.2
$ws.Range("H68").Value = 2085769.9
$ws.Range("I68").Value = 4169034.8
$ws.Range("J68").Value = 2505
$ws.Range("K68").Value = 4169034.8
$ws.Range("L68").Value = 2505
$ws.Range("M68").Value = -4168285.8
$ws.Range("N68").Value = -4003
$ws.Range("H71").Value = 2085769.9
$ws.Range("I71").Value = 4169034.8
$ws.Range("J71").Value = 2505
$ws.Range("K71").Value = 20845174
$ws.Range("L71").Value = 12525
$ws.Range("M71").Value = -20841430
$ws.Range("N71").Value = -20013
$ws.Range("H82").Value = 5641.625
$ws.Range("I82").Value = 3925.8572
$ws.Range("J82").Value = 6976.1113
$ws.Range("K82").Value = 3925.8572
$ws.Range("L82").Value = 6976.1113
$ws.Range("M82").Value = -3564.8572
$ws.Range("N82").Value = -7698.1113
$ws.Range("H85").Value = 5641.625
$ws.Range("I85").Value = 3925.8572
$ws.Range("J85").Value = 6976.1113
$ws.Range("K85").Value = 3925.8572
$ws.Range("L85").Value = 6976.1113
$ws.Range("M85").Value = -2677.8572
$ws.Range("N85").Value = -9472.1113
$ws.Range("H122").Value = 2736.7805
$ws.Range("I122").Value = 2736.7805
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8210.341499999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5760.341499999999
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 6530.304
$ws.Range("I136").Value = 7294.3887
$ws.Range("K136").Value = 21883.1661
$ws.Range("M136").Value = -19333.1661

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8078
$ws.Range("J41").Value = 8078
$ws.Range("L41").Value = 8078
$ws.Range("N41").Value = -8858
$ws.Range("H107").Value = 3230.8333
$ws.Range("I107").Value = 1873.6154
$ws.Range("K107").Value = 5620.8462
$ws.Range("M107").Value = -3700.8462
$ws.Range("H113").Value = 645.17645
$ws.Range("I113").Value = 474.3846
$ws.Range("K113").Value = 1423.1538
$ws.Range("M113").Value = 746.8462
$ws.Range("H132").Value = 304960.53
$ws.Range("J132").Value = 1431189.2
$ws.Range("L132").Value = 4293567.6
$ws.Range("N132").Value = -4298627.6
$ws.Range("H136").Value = 226793.75
$ws.Range("I136").Value = 14726.024
$ws.Range("K136").Value = 44178.072
$ws.Range("M136").Value = -41628.072

Write-Host "All changes applied"